$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1511
$ws.Range("I32").Value = 450
$ws.Range("K32").Value = 450
$ws.Range("M32").Value = -124
$ws.Range("H62").Value = 88177.69500000001
$ws.Range("I62").Value = 141125.5
$ws.Range("J62").Value = 3461.2
$ws.Range("K62").Value = 141125.5
$ws.Range("L62").Value = 3461.2
$ws.Range("M62").Value = -140501.5
$ws.Range("N62").Value = -4709.2
$ws.Range("H65").Value = 88177.69500000001
$ws.Range("I65").Value = 141125.5
$ws.Range("J65").Value = 3461.2
$ws.Range("K65").Value = 705627.5
$ws.Range("L65").Value = 17306
$ws.Range("M65").Value = -702507.5
$ws.Range("N65").Value = -23546
$ws.Range("H129").Value = 2646262.5
$ws.Range("I129").Value = 398.625
$ws.Range("J129").Value = 6174081
$ws.Range("K129").Value = 1195.875
$ws.Range("L129").Value = 18522243
$ws.Range("M129").Value = 3804.125
$ws.Range("N129").Value = -18532243
$ws.Range("H139").Value = 21195.455
$ws.Range("J139").Value = 21195.455
$ws.Range("L139").Value = 21195.455
$ws.Range("N139").Value = -31475.455

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8864.375
$ws.Range("I32").Value = 4789.3384
$ws.Range("K32").Value = 4789.3384
$ws.Range("M32").Value = -4502.3384
$ws.Range("H61").Value = 1073.4546
$ws.Range("I61").Value = 900.8889
$ws.Range("J61").Value = 1850
$ws.Range("K61").Value = 900.8889
$ws.Range("L61").Value = 1850
$ws.Range("M61").Value = -688.8889
$ws.Range("N61").Value = -2274
$ws.Range("H110").Value = 473.5
$ws.Range("I110").Value = 464.70587
$ws.Range("J110").Value = 494.85715
$ws.Range("K110").Value = 464.70587
$ws.Range("L110").Value = 494.85715
$ws.Range("M110").Value = 1580.29413
$ws.Range("N110").Value = -4584.85715
$ws.Range("H122").Value = 1958.5555
$ws.Range("I122").Value = 1919.8
$ws.Range("J122").Value = 2007
$ws.Range("K122").Value = 5759.4
$ws.Range("L122").Value = 6021
$ws.Range("M122").Value = -3309.4
$ws.Range("N122").Value = -10921
$ws.Range("H123").Value = 48796.332
$ws.Range("J123").Value = 48796.332
$ws.Range("L123").Value = 48796.332
$ws.Range("N123").Value = -58596.332
$ws.Range("H136").Value = 1073.4546
$ws.Range("I136").Value = 900.8889
$ws.Range("J136").Value = 1850
$ws.Range("K136").Value = 2702.6667
$ws.Range("L136").Value = 5550
$ws.Range("M136").Value = -152.6667000000002
$ws.Range("N136").Value = -10650

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3570.7585
$ws.Range("I105").Value = 3517.087
$ws.Range("J105").Value = 3776.5
$ws.Range("K105").Value = 3517.087
$ws.Range("L105").Value = 3776.5
$ws.Range("M105").Value = -1770.087
$ws.Range("N105").Value = -7270.5
$ws.Range("H137").Value = 45992.5
$ws.Range("J137").Value = 45992.5
$ws.Range("L137").Value = 45992.5
$ws.Range("N137").Value = -56192.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1835.826
$ws.Range("I31").Value = 1772.5714
$ws.Range("K31").Value = 1772.5714
$ws.Range("M31").Value = -1477.5714
$ws.Range("H34").Value = 1835.826
$ws.Range("I34").Value = 1772.5714
$ws.Range("K34").Value = 1772.5714
$ws.Range("M34").Value = -1570.5714
$ws.Range("H86").Value = 8975.826999999999
$ws.Range("J86").Value = 2166.3845
$ws.Range("L86").Value = 2166.3845
$ws.Range("N86").Value = -4412.3845
$ws.Range("H89").Value = 8975.826999999999
$ws.Range("J89").Value = 2166.3845
$ws.Range("L89").Value = 10831.9225
$ws.Range("N89").Value = -22063.9225
$ws.Range("H99").Value = 2073.2964
$ws.Range("I99").Value = 1546.4117
$ws.Range("K99").Value = 1546.4117
$ws.Range("M99").Value = -48.41170000000011
$ws.Range("H126").Value = 2073.2964
$ws.Range("I126").Value = 1546.4117
$ws.Range("K126").Value = 4639.2351
$ws.Range("M126").Value = -2169.2351

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 794.35
$ws.Range("J131").Value = 893.0964
$ws.Range("L131").Value = 2679.2892
$ws.Range("N131").Value = -12759.2892

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 5845
$ws.Range("I33").Value = 1700
$ws.Range("J33").Value = 9990
$ws.Range("K33").Value = 1700
$ws.Range("L33").Value = 9990
$ws.Range("M33").Value = -1448
$ws.Range("N33").Value = -10494
$ws.Range("H70").Value = 5487.5
$ws.Range("I70").Value = 4316.6665
$ws.Range("J70").Value = 9000
$ws.Range("K70").Value = 4316.6665
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -4046.6665
$ws.Range("N70").Value = -9540
$ws.Range("H73").Value = 5487.5
$ws.Range("I73").Value = 4316.6665
$ws.Range("J73").Value = 9000
$ws.Range("K73").Value = 4316.6665
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -3380.6665
$ws.Range("N73").Value = -10872
$ws.Range("H80").Value = 3180.6296
$ws.Range("I80").Value = 3219.08
$ws.Range("J80").Value = 2700
$ws.Range("K80").Value = 3219.08
$ws.Range("L80").Value = 2700
$ws.Range("M80").Value = -2221.08
$ws.Range("N80").Value = -4696
$ws.Range("H83").Value = 3180.6296
$ws.Range("I83").Value = 3219.08
$ws.Range("J83").Value = 2700
$ws.Range("K83").Value = 16095.4
$ws.Range("L83").Value = 13500
$ws.Range("M83").Value = -11103.4
$ws.Range("N83").Value = -23484
$ws.Range("H122").Value = 880926.25
$ws.Range("I122").Value = 2633978.8
$ws.Range("J122").Value = 4400
$ws.Range("K122").Value = 7901936.399999999
$ws.Range("L122").Value = 13200
$ws.Range("M122").Value = -7899486.399999999
$ws.Range("N122").Value = -18100
$ws.Range("H132").Value = 2536.1562
$ws.Range("I132").Value = 2154.739
$ws.Range("K132").Value = 6464.217000000001
$ws.Range("M132").Value = -3934.217000000001
$ws.Range("H138").Value = 43568.168
$ws.Range("J138").Value = 45281.8
$ws.Range("L138").Value = 45281.8
$ws.Range("N138").Value = -55561.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 13890952
$ws.Range("I100").Value = 18520568
$ws.Range("K100").Value = 18520568
$ws.Range("M100").Value = -18520027
$ws.Range("H122").Value = 3863.4424
$ws.Range("I122").Value = 4888.3335
$ws.Range("J122").Value = 2756.56
$ws.Range("K122").Value = 14665.0005
$ws.Range("L122").Value = 8269.68
$ws.Range("M122").Value = -12215.0005
$ws.Range("N122").Value = -13169.68
$ws.Range("H132").Value = 3920.8823
$ws.Range("I132").Value = 3740.6538
$ws.Range("J132").Value = 4506.625
$ws.Range("K132").Value = 11221.9614
$ws.Range("L132").Value = 13519.875
$ws.Range("M132").Value = -8691.9614
$ws.Range("N132").Value = -18579.875
$ws.Range("H138").Value = 38600
$ws.Range("J138").Value = 38600
$ws.Range("L138").Value = 38600
$ws.Range("N138").Value = -48880

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 29894.8
$ws.Range("J9").Value = 29894.8
$ws.Range("L9").Value = 29894.8
$ws.Range("N9").Value = -30174.8
$ws.Range("H81").Value = 1918.1
$ws.Range("I81").Value = 1772.375
$ws.Range("K81").Value = 3544.75
$ws.Range("M81").Value = -2483.75
$ws.Range("H84").Value = 1918.1
$ws.Range("I84").Value = 1772.375
$ws.Range("K84").Value = 17723.75
$ws.Range("M84").Value = -12419.75
$ws.Range("H107").Value = 862.73914
$ws.Range("I107").Value = 930.0625
$ws.Range("J107").Value = 708.8570999999999
$ws.Range("K107").Value = 2790.1875
$ws.Range("L107").Value = 2126.5713
$ws.Range("M107").Value = -870.1875
$ws.Range("N107").Value = -5966.5713
$ws.Range("H122").Value = 1387.25
$ws.Range("I122").Value = 1333
$ws.Range("K122").Value = 3999
$ws.Range("M122").Value = -1549
$ws.Range("H123").Value = 29378.38
$ws.Range("J123").Value = 29378.38
$ws.Range("L123").Value = 29378.38
$ws.Range("N123").Value = -39178.38
$ws.Range("H132").Value = 1915.4166
$ws.Range("I132").Value = 1934.1765
$ws.Range("J132").Value = 1869.8572
$ws.Range("K132").Value = 5802.529500000001
$ws.Range("L132").Value = 5609.571599999999
$ws.Range("M132").Value = -3272.529500000001
$ws.Range("N132").Value = -10669.5716
$ws.Range("H135").Value = 74679.64
$ws.Range("J135").Value = 74679.64
$ws.Range("L135").Value = 74679.64
$ws.Range("N135").Value = -84819.64
$ws.Range("H136").Value = 1343.6888
$ws.Range("I136").Value = 1391.1082
$ws.Range("J136").Value = 1124.375
$ws.Range("K136").Value = 4173.3246
$ws.Range("L136").Value = 3373.125
$ws.Range("M136").Value = -1623.3246
$ws.Range("N136").Value = -8473.125
$ws.Range("H138").Value = 50946.668
$ws.Range("J138").Value = 50946.668
$ws.Range("L138").Value = 50946.668
$ws.Range("N138").Value = -61226.668
